# Duplicate "Sheet1" twice (mirrors Excel's "Move or Copy… > Create a copy"),
# producing "Sheet1 (2)" and "Sheet1 (3)" right after the original, then make
# "Sheet1 (2)" the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Copy Sheet1 to right after itself -> "Sheet1 (2)"
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# Copy Sheet1 again to right after "Sheet1 (2)" -> "Sheet1 (3)"
$ws1.Copy($null, $ws2)

# Final tab order is now: Sheet1, Sheet1 (2), Sheet1 (3)
# Make "Sheet1 (2)" (index 1, 0-based) the active/visible tab.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
